$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2 value to the new composition text.
$ws.Range("B2").Value = "calcium gluconate (10%w/v) (inj/inf)"

# Add new row 3: A3 stays blank (an empty text cell, matching the existing
# empty A2 cell) and B3 gets the new composition text.
$ws.Range("A3").Value = "'"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "calcium gluconate (10%w/v)"
